$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-32: update Price (D) and/or Volume(1h) (E) ---
$ws.Range("D2").Value = "52.034.74"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "2.934.83"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "357.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "3.395.24"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").Value = "2.932.30"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "52.034.74"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.04%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.188"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +14.27%  "
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.63%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.107"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +14.10%  "
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +0.85%  "

# --- Rows 33-37: ranking shuffled (Toncoin moved to #33); rewrite Coin/Link/Price/Volume ---
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0446"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

# --- Rows 38-51: update Price (D) and/or Volume(1h) (E) ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").Value = "2.141.94"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.247"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.61%  "
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.50%  "
